$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "Table2" structured table by one row (B10:F67 -> B10:F68);
# this keeps the table's ref / autoFilter / worksheet dimension in sync.
$lo = $ws.ListObjects.Item("Table2")
$lo.ListRows.Add() | Out-Null

# New post (#58) data for row 68.
$ws.Range("B68").Value = 58
$ws.Range("C68").Value = "Array variable | Shell Scripting "
$ws.Range("D68").Value = 44178
$ws.Range("E68").Value = "https://programmingport.hashnode.dev/array-variable-or-shell-scripting"
$ws.Range("F68").Value = "https://dev.to/rahulmishra05/array-variable-shell-scripting-56c3"

# Match the formatting used by the rest of the table's data rows.
$ws.Range("B68").Style = $ws.Range("B67").Style
$ws.Range("C68").Style = $ws.Range("C67").Style
$ws.Range("D68").Style = $ws.Range("D67").Style
$ws.Range("D68").NumberFormat = $ws.Range("D67").NumberFormat
$ws.Range("E68").Style = $ws.Range("E67").Style
$ws.Range("F68").Style = $ws.Range("F67").Style

# Match the saved view state: scrolled down with F68 selected.
$ws.Range("F68").Select()
$excel.ActiveWindow.ScrollRow = 51
$excel.ActiveWindow.ScrollColumn = 6
